# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for B2:G7 (TB, d2S, K, IP, Win, sum)
$data = @(
    @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(3.182878228561681, 1.65323645889881,  3.082599426703578,  6.48142807727062,   1, 14.40014219143469),
    @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538),
    @(1.505614041169197, 0.3375848360084654,0.7127328510149897, 0.4998867070740569, 1, 3.055818435266709),
    @(3.182878228561681, 1.65323645889881,  3.082599426703578,  6.48142807727062,   0, 14.40014219143469),
    @(0.006876353814593728, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 1, 1.271808613593681)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
